# Auto update stock data
# The "Date_1" column (A) stores its values as literal text (e.g. "2026/01/04")
# rather than real Excel dates. Every 6th data row (2, 8, 14, ... 74) needs its
# date-looking text bumped from 2026/01/04 to 2026/01/05, while staying plain
# text (Excel would otherwise auto-convert a YYYY/MM/DD-looking string into a
# real date serial when assigned through .Value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "2026/01/04") {
        # Preserve the cell's existing style/number format (plain text)
        # across the write: forcing "@" (Text) right before the assignment
        # stops Excel's auto date-detection from turning the literal
        # "2026/01/05" string into a date serial number.
        $savedStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = "2026/01/05"
        $cell.Style = $savedStyle
    }
}
